$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Allan"
$ws.Range("B1").Value = "Jones"
$ws.Range("C1").Value = "kanalanal"
$ws.Range("D1").Value = 26000

$ws.Range("A2").Value = "Bert"
$ws.Range("B2").Value = "Karlsson"
$ws.Range("C2").Value = "analkanal"
$ws.Range("D2").Value = 3456789

$ws.Range("E1").Value = "User"
$ws.Range("E2").Value = "Admin"

$ws.Range("F1").Value = "Caretaker"

$ws.Range("F1").Select() | Out-Null
